$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 40.397126927018178
$ws.Range("G2").Value = 40.027469031995608
$ws.Range("H2").Value = 40.781982165398048
$ws.Range("I2").Value = 0.00076604310093086566
$ws.Range("J2").Value = 0.00071639483745831306
$ws.Range("K2").Value = 0.0008585975897180203
$ws.Range("L2").Value = 0.057886464662105543
$ws.Range("M2").Value = 0.057526936422573519
$ws.Range("N2").Value = 0.058261002075326598

# Row 3
$ws.Range("F3").Value = 0.0000139073207864358607488
$ws.Range("G3").Value = 0.000000003517178557547588
$ws.Range("H3").Value = 0.000039499859079571012551
$ws.Range("I3").Value = 0.0000121366624099533598259
$ws.Range("J3").Value = 0.00000000327962226382856
$ws.Range("K3").Value = 0.0000343587520125906322583
$ws.Range("L3").Value = 0.000014317086700295460482
$ws.Range("M3").Value = 0.0000000036559564665386888
$ws.Range("N3").Value = 0.0000406425280123063917769

# Row 4
$ws.Range("F4").Value = 40.39714083433897
$ws.Range("G4").Value = 40.027469035512787
$ws.Range("H4").Value = 40.782021665257119
$ws.Range("I4").Value = 0.00077817976334081909
$ws.Range("J4").Value = 0.00071639811708057682
$ws.Range("K4").Value = 0.00089295634173061091
$ws.Range("L4").Value = 0.057900781748805842
$ws.Range("M4").Value = 0.05752694007852998
$ws.Range("N4").Value = 0.058301644603338913
